# Apply "#5: insurance, claim, debt, investment done" changes.
# Sheet "債務" (debt, sheet index 5) and "事業投資" (investment, sheet index 6)
# both get a new "species"/"company" style first column plus seven extra
# trailing metadata columns (property_category .. index), matching the
# layout already used on the other (land/building/deposit/stock) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet5 = 債務 (debt)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# --- header row (row 1) ---
$ws5.Range("B1").Value = "species"
$ws5.Range("C1").Value = "debtor"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "register_date"
$ws5.Range("G1").Value = "register_reason"
$ws5.Range("H1").Value = "property_category"
$ws5.Range("I1").Value = "category"
$ws5.Range("J1").Value = "date"
$ws5.Range("K1").Value = "legislator_name"
$ws5.Range("L1").Value = "legislator_id"
$ws5.Range("M1").Value = "source_file"
$ws5.Range("N1").Value = "index"

# new header cells need the same bold/centered/bordered look as the
# existing header row (B1:G1, style index 1)
$ws5.Range("H1:N1").Font.Bold = $true
$ws5.Range("H1:N1").HorizontalAlignment = -4108
$ws5.Range("H1:N1").VerticalAlignment = -4160
$ws5.Range("H1:N1").Borders.LineStyle = 1

# --- row 2 (legislator property index 128) ---
$ws5.Range("B2").Value = "房屋抵押借款"
$ws5.Range("C2").Value = "呂學樟"
$ws5.Range("D2").Value = "臺灣銀行北大分行"
$ws5.Range("E2").Value = 1242968
$ws5.Range("F2").Value = "93年04月13闩"
$ws5.Range("G2").Value = "購買房屋土地"
$ws5.Range("H2").Value = "debt"
$ws5.Range("I2").Value = "normal"
$ws5.Range("J2").Value = "'2012-04-16"
$ws5.Range("K2").Value = "呂學樟"
$ws5.Range("L2").Value = 892
$ws5.Range("M2").Value = "tmp63271"
$ws5.Range("N2").Value = 128

# --- row 3 (legislator property index 129) ---
$ws5.Range("B3").Value = "土地抵押借款"
$ws5.Range("C3").Value = "呂學樟"
$ws5.Range("D3").Value = "合作金庫商業銀行北新竹分行"
$ws5.Range("E3").Value = 60000000
$ws5.Range("F3").Value = "99年02月08H"
$ws5.Range("G3").Value = "購買土地"
$ws5.Range("H3").Value = "debt"
$ws5.Range("I3").Value = "normal"
$ws5.Range("J3").Value = "'2012-04-16"
$ws5.Range("K3").Value = "呂學樟"
$ws5.Range("L3").Value = 892
$ws5.Range("M3").Value = "tmp63271"
$ws5.Range("N3").Value = 129

# --- row 4 (legislator property index 130) ---
$ws5.Range("B4").Value = "土地抵押借款"
$ws5.Range("C4").Value = "呂學樟"
$ws5.Range("D4").Value = "合作金庫商業銀行北新竹分行"
$ws5.Range("E4").Value = 5000000
$ws5.Range("F4").Value = "100年05月30日"
$ws5.Range("G4").Value = "購買土地"
$ws5.Range("H4").Value = "debt"
$ws5.Range("I4").Value = "normal"
$ws5.Range("J4").Value = "'2012-04-16"
$ws5.Range("K4").Value = "呂學樟"
$ws5.Range("L4").Value = 892
$ws5.Range("M4").Value = "tmp63271"
$ws5.Range("N4").Value = 130

# ---------------------------------------------------------------------
# Sheet6 = 事業投資 (investment)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# --- header row (row 1) ---
$ws6.Range("B1").Value = "owner"
$ws6.Range("C1").Value = "company"
$ws6.Range("D1").Value = "address"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

$ws6.Range("H1:N1").Font.Bold = $true
$ws6.Range("H1:N1").HorizontalAlignment = -4108
$ws6.Range("H1:N1").VerticalAlignment = -4160
$ws6.Range("H1:N1").Borders.LineStyle = 1

# --- row 2 (legislator property index 135) ---
$ws6.Range("B2").Value = "呂學樟"
$ws6.Range("C2").Value = "聯立汽車有限公司"
$ws6.Range("D2").Value = "新竹市中華路一段43號"
$ws6.Range("E2").Value = 13800000
$ws6.Range("F2").Value = "75年03月28日"
$ws6.Range("G2").Value = "成立公司"
$ws6.Range("H2").Value = "investment"
$ws6.Range("I2").Value = "normal"
$ws6.Range("J2").Value = "'2012-04-16"
$ws6.Range("K2").Value = "呂學樟"
$ws6.Range("L2").Value = 892
$ws6.Range("M2").Value = "tmp63271"
$ws6.Range("N2").Value = 135
